$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Column widths below are expressed in the COM "characters" ColumnWidth unit.
# The saved OOXML <col width> ends up pixel-quantized, so these values are
# chosen to land as close as possible to the target widths (~29.98 and
# ~13.75 raw units) after that quantization.
$wideColWidth = 29.166666666666668
$narrowColWidth = 12.833333333333334

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-24 08:51:21"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $narrowColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-24 08:51:28"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $narrowColWidth
